# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Re-order the "Periodo Mora" column (E16:E22) so the most recent period
# (2306) appears first, descending down to the oldest (2212).
$ws.Range("E16").Value = "2306"
$ws.Range("E17").Value = "2305"
$ws.Range("E18").Value = "2304"
$ws.Range("E19").Value = "2303"
$ws.Range("E20").Value = "2302"
$ws.Range("E21").Value = "2301"
$ws.Range("E22").Value = "2212"

# The "Valor Mora" for period 2306 is the special (prorated) value; it now
# lives on row 16, while the previously-first row (2212) takes the regular
# value on row 22.
$ws.Range("F16").Value = 29333
$ws.Range("F22").Value = 40000
